# Insert a new data row at row 81 (pushing the existing rows 81..171 down to
# 82..172) and populate it with the new weekly price-observation record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("81:81").Insert()

$ws.Range("A81").Value = 10
$ws.Range("B81").Value = "Vega Modelo de Temuco"
$ws.Range("C81").Value = "La Araucanía"
$ws.Range("D81").Value = 45079
$ws.Range("E81").Value = 9
$ws.Range("F81").Value = 100112035
$ws.Range("G81").Value = "Bruselas (repollito)"
$ws.Range("H81").Value = "Sin especificar"
$ws.Range("I81").Value = "Primera"
$ws.Range("J81").Value = 25
$ws.Range("K81").Value = 28000
$ws.Range("L81").Value = 28000
$ws.Range("M81").Value = 28000
$ws.Range("N81").Value = "$/malla 15 kilos"
$ws.Range("O81").Value = "Región Metropolitana"
$ws.Range("P81").Value = 1867
$ws.Range("Q81").Value = 15
$ws.Range("R81").Value = "Hortaliza"
